$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version row (row 3): bump the value from 5.0.0 to 6.0.0
$ws.Cells.Item(3, 2).Value = "6.0.0"

# Date row (row 8): update the timestamp
$ws.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher row (row 9): fill in the previously-empty value
$ws.Cells.Item(9, 2).Value = "Alvearie Team"

# Row 10 used to hold "Contact" / "No display for ContactDetail";
# turn it into the "Jurisdiction" / "United States of America" row
$ws.Cells.Item(10, 1).Value = "Jurisdiction"
$ws.Cells.Item(10, 2).Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row - remove it
$ws.Rows.Item(11).Delete()
